$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-77 (columns A:E) reflecting the "burn area effect":
# distance (D) becomes a uniform 20 and travel time (E) becomes a uniform 4
# for every arc, and several arcs (B/C = i/j node pair) were changed/added.
$data = @(
    ,@(1,1,2,20,4)
    ,@(1,1,3,20,4)
    ,@(1,1,5,20,4)
    ,@(1,2,1,20,4)
    ,@(1,2,6,20,4)
    ,@(1,3,1,20,4)
    ,@(1,3,4,20,4)
    ,@(1,3,6,20,4)
    ,@(1,4,3,20,4)
    ,@(1,4,6,20,4)
    ,@(1,4,8,20,4)
    ,@(1,5,1,20,4)
    ,@(1,5,8,20,4)
    ,@(1,6,2,20,4)
    ,@(1,6,3,20,4)
    ,@(1,6,4,20,4)
    ,@(1,6,7,20,4)
    ,@(1,6,9,20,4)
    ,@(1,6,10,20,4)
    ,@(1,7,6,20,4)
    ,@(1,7,8,20,4)
    ,@(1,8,4,20,4)
    ,@(1,8,5,20,4)
    ,@(1,8,7,20,4)
    ,@(1,8,11,20,4)
    ,@(1,9,6,20,4)
    ,@(1,9,10,20,4)
    ,@(1,9,12,20,4)
    ,@(1,9,13,20,4)
    ,@(1,10,6,20,4)
    ,@(1,10,9,20,4)
    ,@(1,10,11,20,4)
    ,@(1,10,13,20,4)
    ,@(1,11,8,20,4)
    ,@(1,11,10,20,4)
    ,@(1,11,14,20,4)
    ,@(1,11,17,20,4)
    ,@(1,12,9,20,4)
    ,@(1,12,15,20,4)
    ,@(1,13,9,20,4)
    ,@(1,13,10,20,4)
    ,@(1,13,15,20,4)
    ,@(1,13,16,20,4)
    ,@(1,14,10,20,4)
    ,@(1,14,11,20,4)
    ,@(1,14,16,20,4)
    ,@(1,14,17,20,4)
    ,@(1,15,12,20,4)
    ,@(1,15,13,20,4)
    ,@(1,15,18,20,4)
    ,@(1,16,13,20,4)
    ,@(1,16,14,20,4)
    ,@(1,16,19,20,4)
    ,@(1,16,21,20,4)
    ,@(1,17,11,20,4)
    ,@(1,17,14,20,4)
    ,@(1,17,19,20,4)
    ,@(1,17,20,20,4)
    ,@(1,18,15,20,4)
    ,@(1,18,21,20,4)
    ,@(1,18,23,20,4)
    ,@(1,19,16,20,4)
    ,@(1,19,17,20,4)
    ,@(1,19,21,20,4)
    ,@(1,20,17,20,4)
    ,@(1,20,22,20,4)
    ,@(1,21,13,20,4)
    ,@(1,21,16,20,4)
    ,@(1,21,18,20,4)
    ,@(1,21,19,20,4)
    ,@(1,21,22,20,4)
    ,@(1,21,23,20,4)
    ,@(1,22,20,20,4)
    ,@(1,22,21,20,4)
    ,@(1,23,18,20,4)
    ,@(1,23,21,20,4)
)

$startRow = 2
for ($idx = 0; $idx -lt $data.Count; $idx++) {
    $r = $startRow + $idx
    $rowVals = $data[$idx]
    $ws.Cells.Item($r, 1).Value = $rowVals[0]
    $ws.Cells.Item($r, 2).Value = $rowVals[1]
    $ws.Cells.Item($r, 3).Value = $rowVals[2]
    $ws.Cells.Item($r, 4).Value = $rowVals[3]
    $ws.Cells.Item($r, 5).Value = $rowVals[4]
}

# Update the view: scrolled down with a new active selection, matching the
# author's final cursor position after adding the extra burn-area rows.
$ws.Application.ActiveWindow.ScrollRow = 54
$ws.Range("H69").Select()
